$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A holds "dates" stored as plain text (inline strings), each
# shifted forward by one week (7 days). Force text format first so
# Excel does not auto-convert the strings into date serial numbers.
$dates = @{
    1  = "2020-03-12"
    2  = "2020-03-12"
    3  = "2020-03-08"
    4  = "2020-03-12"
    5  = "2020-03-15"
    6  = "2020-03-19"
    7  = "2020-03-22"
    8  = "2020-03-26"
    9  = "2020-03-29"
    10 = "2020-04-02"
    11 = "2020-04-05"
    12 = "2020-04-09"
    13 = "2020-04-12"
    14 = "2020-04-16"
    15 = "2020-04-12"
}

foreach ($row in $dates.Keys) {
    $cell = $ws.Range("A$row")
    $cell.NumberFormat = "@"
    $cell.Value = $dates[$row]
}
